$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 221
$ws.Range("F5").Value = 8984
$ws.Range("F6").Value = 521
$ws.Range("F7").Value = 89
$ws.Range("F9").Value = 178
$ws.Range("F10").Value = 319
$ws.Range("F11").Value = 365
$ws.Range("F15").Value = 396
$ws.Range("F16").Value = 11455
$ws.Range("F18").Value = 303
$ws.Range("F24").Value = 209
$ws.Range("F26").Value = 166
$ws.Range("F27").Value = 97
$ws.Range("F29").Value = 2700
$ws.Range("F31").Value = 99
$ws.Range("F32").Value = 2087
$ws.Range("F33").Value = 49
$ws.Range("F34").Value = 2123
$ws.Range("F35").Value = 941
$ws.Range("F36").Value = 4131
$ws.Range("F38").Value = 302
$ws.Range("F39").Value = 2598
$ws.Range("F40").Value = 3041
$ws.Range("F41").Value = 1273
$ws.Range("F44").Value = 370
$ws.Range("F45").Value = 400
$ws.Range("F46").Value = 58
$ws.Range("F47").Value = 157
$ws.Range("F48").Value = 99

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 17
$ws.Range("F6").Value = 9
$ws.Range("F19").Value = 34

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 41

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 17
$ws.Range("F5").Value = 9
$ws.Range("F7").Value = 221
$ws.Range("F9").Value = 8984
$ws.Range("F10").Value = 521
$ws.Range("F13").Value = 178
$ws.Range("F14").Value = 319
$ws.Range("F15").Value = 365
$ws.Range("F18").Value = 11455
$ws.Range("F19").Value = 303
$ws.Range("F23").Value = 41
$ws.Range("F28").Value = 97
$ws.Range("F30").Value = 2700
$ws.Range("F31").Value = 99
$ws.Range("F32").Value = 2087
$ws.Range("F33").Value = 49
$ws.Range("F34").Value = 941
$ws.Range("F35").Value = 4131
$ws.Range("F37").Value = 302
$ws.Range("F38").Value = 2598
$ws.Range("F39").Value = 3041
$ws.Range("F41").Value = 1273
$ws.Range("F43").Value = 370
$ws.Range("F45").Value = 400
$ws.Range("F46").Value = 58
$ws.Range("F47").Value = 157
$ws.Range("F48").Value = 99
